$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 64 new rows at the top, pushing the existing "Events" (9.2) checklist block
# down from rows 1-32 to rows 69-96.
$ws.Rows("1:64").Insert()

# Fill in the new rows with the two additional checklist blocks:
# "Digital Content-management" (9.7) in rows 1-32
# "POS Communication" (9.3) in rows 37-64
$ws.Range("A5").Value = "1: must`n(if POS content management tool is available in the market)"
$ws.Range("A7").Value = 0
$ws.Range("A9").Value = "N"
$ws.Range("A11").Value = "Marketing Manager"
$ws.Range("A13").Value = "> visual check "
$ws.Range("A15").Value = "Is the most actual content is shown on all digital devices? (if not used by customers or sales personnel for individual product presentation)"
$ws.Range("A17").Value = "X"
$ws.Range("A19").Value = "Audit"
$ws.Range("A21").Value = "X"
$ws.Range("A22").Value = "X"
$ws.Range("A23").Value = "X"
$ws.Range("A24").Value = "X"
$ws.Range("A25").Value = "X"
$ws.Range("A26").Value = "X"
$ws.Range("A28").Value = " (if POS content management tool is available in the market)"
$ws.Range("A30").Value = "The importer ensures that only the most actual content is shown on all digital devices.`n"
$ws.Range("A31").Formula = '="9.7"'
$ws.Range("A31").Copy() | Out-Null
$ws.Range("A31").PasteSpecial(-4163) | Out-Null
$ws.Range("A32").Value = "Digital Content-management"
$ws.Range("A37").Value = "1: must `n(if Regional Office offers POS Communication)"
$ws.Range("A39").Value = 0
$ws.Range("A41").Value = "N"
$ws.Range("A43").Value = "Marketing Manager"
$ws.Range("A45").Value = "> check of orders/deliveries "
$ws.Range("A47").Value = "Is the material provided for launches (POS) ordered in time?"
$ws.Range("A49").Value = "C"
$ws.Range("A51").Value = "Regional Office or Importer"
$ws.Range("A53").Value = "X"
$ws.Range("A54").Value = "X"
$ws.Range("A55").Value = "X"
$ws.Range("A56").Value = "X"
$ws.Range("A57").Value = "X"
$ws.Range("A58").Value = "X"
$ws.Range("A62").Value = "The importer orders the material provided for launches (Point of Sale (POS) Launch Kit) and other material for further product/technology or seasonal related topics in time and implements this in accordance with BMW Regional Office specifications.`n"
$ws.Range("A63").Formula = '="9.3"'
$ws.Range("A63").Copy() | Out-Null
$ws.Range("A63").PasteSpecial(-4163) | Out-Null
$ws.Range("A64").Value = "POS Communication"
